$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.561.01"
$ws.Range("E2").Value = "  +1.56%  "

$ws.Range("D3").Value = "1.911.01"
$ws.Range("E3").Value = "  +3.41%  "

$ws.Range("E4").Value = "  +0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.90"
$ws.Range("E5").Value = "  +5.23%  "

$ws.Range("E6").Value = "  +1.68%  "

$ws.Range("E7").Value = "  +0.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.18"
$ws.Range("E8").Value = "  -0.55%  "

$ws.Range("E9").Value = "  +2.99%  "

$ws.Range("E10").Value = "  +1.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0998"
$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("D12").Value = "2.184.42"
$ws.Range("E12").Value = "  +3.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.46"
$ws.Range("E13").Value = "  +9.86%  "

$ws.Range("D14").Value = "1.924.67"
$ws.Range("E14").Value = "  +3.08%  "

$ws.Range("E15").Value = "  +2.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.88"
$ws.Range("E16").Value = "  +4.13%  "

$ws.Range("D17").Value = "35.514.44"
$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "72.19"
$ws.Range("E18").Value = "  +3.21%  "

$ws.Range("D19").Value = "0.0₃0822"
$ws.Range("E19").Value = "  +3.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.83"
$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("E21").Value = "  +3.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.90"
$ws.Range("E22").Value = "  +2.23%  "

$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.31"
$ws.Range("E25").Value = "  +0.81%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  +20.33%  "

$ws.Range("E27").Value = "  +9.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.03"
$ws.Range("E28").Value = "  +2.50%  "

$ws.Range("E29").Value = "  +0.83%  "

$ws.Range("E30").Value = "  +27.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0571"
$ws.Range("E31").Value = "  +2.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.12"
$ws.Range("E32").Value = "  +3.78%  "

$ws.Range("E33").Value = "  +5.98%  "

$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.75"
$ws.Range("E35").Value = "  +6.34%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.40"
$ws.Range("E36").Value = "  +13.36%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.03"
$ws.Range("E37").Value = "  +2.37%  "

$ws.Range("E38").Value = "  +3.74%  "

$ws.Range("E39").Value = "  +1.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "91.26"
$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.71"
$ws.Range("E41").Value = "  +5.51%  "

$ws.Range("D42").Value = "1.356.91"
$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.43"
$ws.Range("E43").Value = "  +42.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0603"
$ws.Range("E44").Value = "  +13.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").Value = "  +2.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.68"
$ws.Range("E46").Value = "  +1.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +1.38%  "

$ws.Range("E48").Value = "  +0.51%  "

$ws.Range("E49").Value = "  +4.99%  "

$ws.Range("D50").Value = "2.093.42"
$ws.Range("E50").Value = "  +3.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0692"
$ws.Range("E51").Value = "  +2.00%  "
